# PowerShell COM-interop script (iron_native / Excel.Application)
# Re-applies a live crypto-market data refresh captured 2024-11-22 06:01:32:
#   - "Top 50 Cryptocurrencies": updated price / market-cap / volume / 24h-change
#     for (almost) every row, plus Kaspa/Render swapping rank positions 37 <-> 38
#   - "Top 5 by Market Cap": refreshed market-cap figures to stay in sync
#   - "Summary": refreshed the three derived headline metrics
#
# Numbers are written in plain decimal (never scientific notation) because the
# COM-interop PowerShell parser here does not accept tokens like "5.213e-05".

$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet: "Top 50 Cryptocurrencies"
# ======================================================================
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

# Rows 37 & 38 swap places: Render <-> Kaspa (name + symbol)
$ws1.Cells.Item(37, 1).Value = "Kaspa"
$ws1.Cells.Item(37, 2).Value = "kas"
$ws1.Cells.Item(38, 1).Value = "Render"
$ws1.Cells.Item(38, 2).Value = "render"

# Row 2: C, D, E, F column(s) changed
$ws1.Cells.Item(2, 3).Value = 98950
$ws1.Cells.Item(2, 4).Value = 1957801914152
$ws1.Cells.Item(2, 5).Value = 99838698778
$ws1.Cells.Item(2, 6).Value = 2.0972

# Row 3: C, D, E, F column(s) changed
$ws1.Cells.Item(3, 3).Value = 3389.14
$ws1.Cells.Item(3, 4).Value = 408183778081
$ws1.Cells.Item(3, 5).Value = 57101803551
$ws1.Cells.Item(3, 6).Value = 9.253769999999999

# Row 4: D, E, F column(s) changed
$ws1.Cells.Item(4, 4).Value = 130823802841
$ws1.Cells.Item(4, 5).Value = 191268654241
$ws1.Cells.Item(4, 6).Value = 0.32316

# Row 5: C, D, E, F column(s) changed
$ws1.Cells.Item(5, 3).Value = 261.36
$ws1.Cells.Item(5, 4).Value = 124099389605
$ws1.Cells.Item(5, 5).Value = 14903347184
$ws1.Cells.Item(5, 6).Value = 9.29208

# Row 6: C, D, E, F column(s) changed
$ws1.Cells.Item(6, 3).Value = 633.61
$ws1.Cells.Item(6, 4).Value = 92441290124
$ws1.Cells.Item(6, 5).Value = 2462097536
$ws1.Cells.Item(6, 6).Value = 4.16909

# Row 7: D, E, F column(s) changed
$ws1.Cells.Item(7, 4).Value = 79890440378
$ws1.Cells.Item(7, 5).Value = 18216983753
$ws1.Cells.Item(7, 6).Value = 26.85773

# Row 8: C, D, E, F column(s) changed
$ws1.Cells.Item(8, 3).Value = 0.396715
$ws1.Cells.Item(8, 4).Value = 58265317109
$ws1.Cells.Item(8, 5).Value = 9787067581
$ws1.Cells.Item(8, 6).Value = 3.75439

# Row 9: C, D, E, F column(s) changed
$ws1.Cells.Item(9, 3).Value = 0.999911
$ws1.Cells.Item(9, 4).Value = 38332070800
$ws1.Cells.Item(9, 5).Value = 8044180063
$ws1.Cells.Item(9, 6).Value = 0.38615

# Row 10: C, D, E, F column(s) changed
$ws1.Cells.Item(10, 3).Value = 3387.07
$ws1.Cells.Item(10, 4).Value = 33211105579
$ws1.Cells.Item(10, 5).Value = 145751955
$ws1.Cells.Item(10, 6).Value = 9.066129999999999

# Row 11: C, D, E, F column(s) changed
$ws1.Cells.Item(11, 3).Value = 0.882131
$ws1.Cells.Item(11, 4).Value = 31574566650
$ws1.Cells.Item(11, 5).Value = 3012421450
$ws1.Cells.Item(11, 6).Value = 12.58447

# Row 12: C, D, E, F column(s) changed
$ws1.Cells.Item(12, 3).Value = 0.200153
$ws1.Cells.Item(12, 4).Value = 17285755290
$ws1.Cells.Item(12, 5).Value = 1062451930
$ws1.Cells.Item(12, 6).Value = 1.88472

# Row 13: C, D, E, F column(s) changed
$ws1.Cells.Item(13, 3).Value = 36.38
$ws1.Cells.Item(13, 4).Value = 14894142924
$ws1.Cells.Item(13, 5).Value = 1049098992
$ws1.Cells.Item(13, 6).Value = 7.78136

# Row 14: D, E, F column(s) changed
$ws1.Cells.Item(14, 4).Value = 14745167792
$ws1.Cells.Item(14, 5).Value = 1599456205
$ws1.Cells.Item(14, 6).Value = 4.18302

# Row 15: C, D, E, F column(s) changed
$ws1.Cells.Item(15, 3).Value = 4011.1
$ws1.Cells.Item(15, 4).Value = 14487010535
$ws1.Cells.Item(15, 5).Value = 169572218
$ws1.Cells.Item(15, 6).Value = 9.07639

# Row 16: C, D, E, F column(s) changed
$ws1.Cells.Item(16, 3).Value = 98859
$ws1.Cells.Item(16, 4).Value = 14444416042
$ws1.Cells.Item(16, 5).Value = 847328202
$ws1.Cells.Item(16, 6).Value = 2.29471

# Row 17: D, E, F column(s) changed
$ws1.Cells.Item(17, 4).Value = 14159481773
$ws1.Cells.Item(17, 5).Value = 637235491
$ws1.Cells.Item(17, 6).Value = 4.09133

# Row 18: D, E, F column(s) changed
$ws1.Cells.Item(18, 4).Value = 10273589779
$ws1.Cells.Item(18, 5).Value = 2444539501
$ws1.Cells.Item(18, 6).Value = 1.08414

# Row 19: C, D, E, F column(s) changed
$ws1.Cells.Item(19, 3).Value = 498.78
$ws1.Cells.Item(19, 4).Value = 9871422041
$ws1.Cells.Item(19, 5).Value = 1867669705
$ws1.Cells.Item(19, 6).Value = -2.15394

# Row 20: C, D, E, F column(s) changed
$ws1.Cells.Item(20, 3).Value = 3392
$ws1.Cells.Item(20, 4).Value = 9734929278
$ws1.Cells.Item(20, 5).Value = 2223032283
$ws1.Cells.Item(20, 6).Value = 9.49335

# Row 21: D, E, F column(s) changed
$ws1.Cells.Item(21, 4).Value = 9592025550
$ws1.Cells.Item(21, 5).Value = 1257405474
$ws1.Cells.Item(21, 6).Value = 5.48991

# Row 22: D, E, F column(s) changed
$ws1.Cells.Item(22, 4).Value = 8996309587
$ws1.Cells.Item(22, 5).Value = 6804659980
$ws1.Cells.Item(22, 6).Value = 10.57673

# Row 23: C, D, E, F column(s) changed
$ws1.Cells.Item(23, 3).Value = 6.22
$ws1.Cells.Item(23, 4).Value = 8965227660
$ws1.Cells.Item(23, 5).Value = 830110510
$ws1.Cells.Item(23, 6).Value = 10.00362

# Row 24: C, D, E, F column(s) changed
$ws1.Cells.Item(24, 3).Value = 0.285028
$ws1.Cells.Item(24, 4).Value = 8544502079
$ws1.Cells.Item(24, 5).Value = 2312821480
$ws1.Cells.Item(24, 6).Value = 21.06954

# Row 25: C, D, E, F column(s) changed
$ws1.Cells.Item(25, 3).Value = 8.81
$ws1.Cells.Item(25, 4).Value = 8122298571
$ws1.Cells.Item(25, 5).Value = 3437792
$ws1.Cells.Item(25, 6).Value = 4.04549

# Row 26: C, D, E, F column(s) changed
$ws1.Cells.Item(26, 3).Value = 5.82
$ws1.Cells.Item(26, 4).Value = 7088167865
$ws1.Cells.Item(26, 5).Value = 1012128839
$ws1.Cells.Item(26, 6).Value = 5.59196

# Row 27: C, D, E, F column(s) changed
$ws1.Cells.Item(27, 3).Value = 90.94
$ws1.Cells.Item(27, 4).Value = 6841601458
$ws1.Cells.Item(27, 5).Value = 1410008469
$ws1.Cells.Item(27, 6).Value = 5.26465

# Row 28: C, D, E, F column(s) changed
$ws1.Cells.Item(28, 3).Value = 12.14
$ws1.Cells.Item(28, 4).Value = 6471544641
$ws1.Cells.Item(28, 5).Value = 859704634
$ws1.Cells.Item(28, 6).Value = 4.33772

# Row 29: C, D, E, F column(s) changed
$ws1.Cells.Item(29, 3).Value = 3576.2
$ws1.Cells.Item(29, 4).Value = 6220401781
$ws1.Cells.Item(29, 5).Value = 106018794
$ws1.Cells.Item(29, 6).Value = 9.409129999999999

# Row 30: C, D, E, F column(s) changed
$ws1.Cells.Item(30, 3).Value = 9.41
$ws1.Cells.Item(30, 4).Value = 5655338943
$ws1.Cells.Item(30, 5).Value = 866236092
$ws1.Cells.Item(30, 6).Value = 6.92099

# Row 31: C, D, E, F column(s) changed
$ws1.Cells.Item(31, 3).Value = 0.203305
$ws1.Cells.Item(31, 4).Value = 5504574071
$ws1.Cells.Item(31, 5).Value = 132475693
$ws1.Cells.Item(31, 6).Value = 16.57485

# Row 32: C, D, E, F column(s) changed
$ws1.Cells.Item(32, 3).Value = 0.99746
$ws1.Cells.Item(32, 4).Value = 5229388035
$ws1.Cells.Item(32, 5).Value = 91130
$ws1.Cells.Item(32, 6).Value = 0.04723

# Row 33: C, D, E, F column(s) changed
$ws1.Cells.Item(33, 3).Value = 0.133813
$ws1.Cells.Item(33, 4).Value = 5110910421
$ws1.Cells.Item(33, 5).Value = 907416623
$ws1.Cells.Item(33, 6).Value = 8.14057

# Row 34: C, D, E, F column(s) changed
$ws1.Cells.Item(34, 3).Value = 9.640000000000001
$ws1.Cells.Item(34, 4).Value = 4572289134
$ws1.Cells.Item(34, 5).Value = 274975740
$ws1.Cells.Item(34, 6).Value = 7.32213

# Row 35: C, D, E, F column(s) changed
$ws1.Cells.Item(35, 3).Value = 28.1
$ws1.Cells.Item(35, 4).Value = 4205261170
$ws1.Cells.Item(35, 5).Value = 871015505
$ws1.Cells.Item(35, 6).Value = 6.17713

# Row 36: C, D, E, F column(s) changed
$ws1.Cells.Item(36, 3).Value = 0.00005213
$ws1.Cells.Item(36, 4).Value = 3915257516
$ws1.Cells.Item(36, 5).Value = 1670661215
$ws1.Cells.Item(36, 6).Value = 2.57426

# Row 37: C, D, E, F column(s) changed
$ws1.Cells.Item(37, 3).Value = 0.152234
$ws1.Cells.Item(37, 4).Value = 3835490531
$ws1.Cells.Item(37, 5).Value = 151673221
$ws1.Cells.Item(37, 6).Value = 1.06944

# Row 38: C, D, E, F column(s) changed
$ws1.Cells.Item(38, 3).Value = 7.41
$ws1.Cells.Item(38, 4).Value = 3835014471
$ws1.Cells.Item(38, 5).Value = 432564798
$ws1.Cells.Item(38, 6).Value = 1.19444

# Row 39: C, D, E, F column(s) changed
$ws1.Cells.Item(39, 3).Value = 0.473232
$ws1.Cells.Item(39, 4).Value = 3772214715
$ws1.Cells.Item(39, 5).Value = 476956609
$ws1.Cells.Item(39, 6).Value = 8.94638

# Row 40: C, D, E, F column(s) changed
$ws1.Cells.Item(40, 3).Value = 504.62
$ws1.Cells.Item(40, 4).Value = 3728030742
$ws1.Cells.Item(40, 5).Value = 282237842
$ws1.Cells.Item(40, 6).Value = 2.92676

# Row 41: D, E, F column(s) changed
$ws1.Cells.Item(41, 4).Value = 3688538085
$ws1.Cells.Item(41, 5).Value = 225145685
$ws1.Cells.Item(41, 6).Value = -0.05451

# Row 42: D, E, F column(s) changed
$ws1.Cells.Item(42, 4).Value = 3575035917
$ws1.Cells.Item(42, 5).Value = 33355617
$ws1.Cells.Item(42, 6).Value = 3.25778

# Row 43: C, D, E, F column(s) changed
$ws1.Cells.Item(43, 3).Value = 0.999623
$ws1.Cells.Item(43, 4).Value = 3440873426
$ws1.Cells.Item(43, 5).Value = 153373746
$ws1.Cells.Item(43, 6).Value = 0.32737

# Row 44: C, D, E, F column(s) changed
$ws1.Cells.Item(44, 3).Value = 3.39
$ws1.Cells.Item(44, 4).Value = 3388980858
$ws1.Cells.Item(44, 5).Value = 1282763387
$ws1.Cells.Item(44, 6).Value = 6.66468

# Row 45: C, D, E, F column(s) changed
$ws1.Cells.Item(45, 3).Value = 3.73
$ws1.Cells.Item(45, 4).Value = 3363633788
$ws1.Cells.Item(45, 5).Value = 301881074
$ws1.Cells.Item(45, 6).Value = 4.29752

# Row 46: D, E, F column(s) changed
$ws1.Cells.Item(46, 4).Value = 3347597945
$ws1.Cells.Item(46, 5).Value = 485036679
$ws1.Cells.Item(46, 6).Value = 3.18392

# Row 47: C, D, E, F column(s) changed
$ws1.Cells.Item(47, 3).Value = 0.792019
$ws1.Cells.Item(47, 4).Value = 3245886344
$ws1.Cells.Item(47, 5).Value = 1672660052
$ws1.Cells.Item(47, 6).Value = 14.82186

# Row 48: C, D, E, F column(s) changed
$ws1.Cells.Item(48, 3).Value = 161.12
$ws1.Cells.Item(48, 4).Value = 2972253972
$ws1.Cells.Item(48, 5).Value = 86424934
$ws1.Cells.Item(48, 6).Value = -0.18419

# Row 49: C, D, E, F column(s) changed
$ws1.Cells.Item(49, 3).Value = 1.97
$ws1.Cells.Item(49, 4).Value = 2958504435
$ws1.Cells.Item(49, 5).Value = 358440221
$ws1.Cells.Item(49, 6).Value = 3.285

# Row 50: C, D, E, F column(s) changed
$ws1.Cells.Item(50, 3).Value = 0.8428
$ws1.Cells.Item(50, 4).Value = 2837535662
$ws1.Cells.Item(50, 5).Value = 184150203
$ws1.Cells.Item(50, 6).Value = 16.17144

# Row 51: C, D, E, F column(s) changed
$ws1.Cells.Item(51, 3).Value = 4.71
$ws1.Cells.Item(51, 4).Value = 2830209760
$ws1.Cells.Item(51, 5).Value = 578155150
$ws1.Cells.Item(51, 6).Value = 8.21808

# ======================================================================
# Sheet: "Top 5 by Market Cap"
# ======================================================================
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws2.Range("B2").Value = 1957801914152
$ws2.Range("B3").Value = 408183778081
$ws2.Range("B4").Value = 130823802841
$ws2.Range("B5").Value = 124099389605
$ws2.Range("B6").Value = 92441290124

# ======================================================================
# Sheet: "Summary"
# ======================================================================
$ws3 = $wb.Worksheets.Item("Summary")

# B2 starts with "$" which Excel/COM auto-coerces a plain .Value assignment
# into a currency NUMBER (losing the literal "$" text). Force text mode for
# the assignment, then restore the cell to the default ("Normal") style so no
# stray number-format is left behind on the cell.
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "$4358.28"
$ws3.Range("B2").Style = "Normal"

$ws3.Range("B3").Value = "XRP (26.86%)"
$ws3.Range("B4").Value = "Bitcoin Cash (-2.15%)"
